$wb = $excel.ActiveWorkbook

$newMdName = "cea20f97-1964-4eb3-9c5e-97c025d67b6eooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newMdPath = "e2e\" + $newMdName
$newZhXlf = "cea20f97-1964-4eb3-9c5e-97c025d67b6eoooooooooooooooooooooooooooooooooooooooo.9f277ebd71f7c7bd24f62d659a8b3bd0f94e6c3f.zh-cn.xlf"
$newDeXlf = "cea20f97-1964-4eb3-9c5e-97c025d67b6eoooooooooooooooooooooooooooooooooooooooo.9f277ebd71f7c7bd24f62d659a8b3bd0f94e6c3f.de-de.xlf"
$newHoDate = "2016-08-31 04:28:38"
$newZhDate = "2016-08-31 04:28:33"
$newDeDate = "2016-08-31 04:28:38"
$status = "Ready for handoff"
$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/97e6375357c64a04e4ad7b2006e72161d03ddbeb/e2e/" + $newMdName

# ---------- Overview sheet ----------
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)
$rowOv = $loOv.ListRows.Add()
$wsOv.Range("A3").Value = $newMdName
$wsOv.Range("B3").Value = $newMdPath
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("D3").Value = ""
$wsOv.Range("E3").Value = $status
$wsOv.Range("F3").Value = $status
$wsOv.Range("G3").Value = $newHoDate
$wsOv.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), $hyperlinkUrl, [Type]::Missing, [Type]::Missing, $newMdPath) | Out-Null

# ---------- zh-cn sheet ----------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$rowZh = $loZh.ListRows.Add()
$wsZh.Range("A3").Value = $newMdName
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = $status
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = $newZhXlf
$wsZh.Range("H3").Value = $newZhDate
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $hyperlinkUrl, [Type]::Missing, [Type]::Missing, $newMdName) | Out-Null

# ---------- de-de sheet ----------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$rowDe = $loDe.ListRows.Add()
$wsDe.Range("A3").Value = $newMdName
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = $status
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = $newDeXlf
$wsDe.Range("H3").Value = $newDeDate
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $hyperlinkUrl, [Type]::Missing, [Type]::Missing, $newMdName) | Out-Null
